# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity, and Temperature sheets,
# extending each sheet's used range. All columns in these logs are stored
# as plain text (dates/times/percentages as literal strings, not native
# Excel types), so we force the destination cells to Text format ("@")
# before writing values -- this stops Excel from auto-converting things
# like "2026-02-06" into a date serial or "68.1%" into a numeric percentage.

$wb = $excel.ActiveWorkbook

# NOTE: call this with POSITIONAL args (not -Worksheet/-Rows). Passing an
# array via a named parameter to a function in this interpreter loses the
# array contents (ends up with Count 0), while positional passing works.
function Append-LogRows {
    param($Worksheet, $Rows)

    $firstRow = $Rows[0][0]
    $lastRow = $Rows[$Rows.Count - 1][0]

    # Pre-format the destination block as Text so values stick verbatim.
    $destRange = $Worksheet.Range($Worksheet.Cells.Item($firstRow, 1), $Worksheet.Cells.Item($lastRow, 6))
    $destRange.NumberFormat = "@"

    foreach ($row in $Rows) {
        $r = $row[0]
        $Worksheet.Cells.Item($r, 1).Value = $row[1]
        $Worksheet.Cells.Item($r, 2).Value = $row[2]
        $Worksheet.Cells.Item($r, 3).Value = $row[3]
        $Worksheet.Cells.Item($r, 4).Value = $row[4]
        $Worksheet.Cells.Item($r, 5).Value = $row[5]
        $Worksheet.Cells.Item($r, 6).Value = $row[6]
    }
}

# --- PIR sheet: append rows 500-512 -----------------------------------------
$wsPIR = $wb.Worksheets.Item("PIR")
$pirRows = @(
    @(500, "2026-02-06", "10:21:30", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(501, "2026-02-06", "10:21:33", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(502, "2026-02-06", "10:21:37", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(503, "2026-02-06", "10:21:40", "10:00", "Bathroom", "Motion Detected", "Active"),
    @(504, "2026-02-06", "10:21:44", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(505, "2026-02-06", "10:21:49", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(506, "2026-02-06", "10:21:54", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(507, "2026-02-06", "10:21:59", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(508, "2026-02-06", "10:22:04", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(509, "2026-02-06", "10:22:10", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(510, "2026-02-06", "10:22:14", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(511, "2026-02-06", "10:22:20", "10:00", "Bathroom", "No Motion", "Inactive"),
    @(512, "2026-02-06", "10:22:25", "10:00", "Bathroom", "No Motion", "Inactive")
)
Append-LogRows $wsPIR $pirRows

# --- Humidity sheet: append rows 350-359 ------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRows = @(
    @(350, "2026-02-06", "10:21:28", "10:00", "Bathroom", "68.1%", "Active"),
    @(351, "2026-02-06", "10:21:31", "10:00", "Bathroom", "67.2%", "Active"),
    @(352, "2026-02-06", "10:21:34", "10:00", "Bathroom", "68.2%", "Active"),
    @(353, "2026-02-06", "10:21:38", "10:00", "Bathroom", "67.3%", "Active"),
    @(354, "2026-02-06", "10:21:41", "10:00", "Bathroom", "68.3%", "Active"),
    @(355, "2026-02-06", "10:21:46", "10:00", "Bathroom", "67.4%", "Active"),
    @(356, "2026-02-06", "10:21:51", "10:00", "Bathroom", "68.4%", "Active"),
    @(357, "2026-02-06", "10:21:56", "10:00", "Bathroom", "67.4%", "Active"),
    @(358, "2026-02-06", "10:22:01", "10:00", "Bathroom", "68.4%", "Active"),
    @(359, "2026-02-06", "10:22:11", "10:00", "Bathroom", "68.4%", "Active")
)
Append-LogRows $wsHumidity $humidityRows

# --- Temperature sheet: append rows 350-359 ---------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRows = @(
    @(350, "2026-02-06", "10:21:29", "10:00", "Bathroom", "28.2C", "Active"),
    @(351, "2026-02-06", "10:21:32", "10:00", "Bathroom", "28.2C", "Active"),
    @(352, "2026-02-06", "10:21:36", "10:00", "Bathroom", "28.2C", "Active"),
    @(353, "2026-02-06", "10:21:39", "10:00", "Bathroom", "28.2C", "Active"),
    @(354, "2026-02-06", "10:21:42", "10:00", "Bathroom", "28.2C", "Active"),
    @(355, "2026-02-06", "10:21:47", "10:00", "Bathroom", "28.2C", "Active"),
    @(356, "2026-02-06", "10:21:52", "10:00", "Bathroom", "28.2C", "Active"),
    @(357, "2026-02-06", "10:21:57", "10:00", "Bathroom", "28.2C", "Active"),
    @(358, "2026-02-06", "10:22:02", "10:00", "Bathroom", "28.2C", "Active"),
    @(359, "2026-02-06", "10:22:12", "10:00", "Bathroom", "28.2C", "Active")
)
Append-LogRows $wsTemperature $temperatureRows

Write-Output "Appended $($pirRows.Count) rows to PIR, $($humidityRows.Count) rows to Humidity, $($temperatureRows.Count) rows to Temperature."
